$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textRows = @(5, 6, 10, 11, 12, 20, 21, 23, 24, 26, 28, 31, 33, 34, 37, 38, 40, 41, 42, 45, 47, 50, 51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = '49.718.23'
$ws.Range("E2").Value = '  +3.16%  '

$ws.Range("D3").Value = '2.616.16'
$ws.Range("E3").Value = '  +4.43%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '324.48'
$ws.Range("E5").Value = '  +0.97%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '110.06'
$ws.Range("E6").Value = '  +1.65%  '

$ws.Range("E7").Value = '  +1.49%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  +3.70%  '

$ws.Range("D10").Value = '40.87'
$ws.Range("E10").Value = '  +2.28%  '

$ws.Range("D11").Value = '20.69'
$ws.Range("E11").Value = '  +2.23%  '

$ws.Range("D12").Value = '0.0822'
$ws.Range("E12").Value = '  +0.51%  '

$ws.Range("E13").Value = '  +0.79%  '

$ws.Range("E14").Value = '  +1.85%  '

$ws.Range("D15").Value = '3.026.40'
$ws.Range("E15").Value = '  +4.58%  '

$ws.Range("D16").Value = '2.611.35'
$ws.Range("E16").Value = '  +4.19%  '

$ws.Range("E17").Value = '  +3.30%  '

$ws.Range("D18").Value = '49.682.31'
$ws.Range("E18").Value = '  +3.42%  '

$ws.Range("E19").Value = '  +11.46%  '

$ws.Range("D20").Value = '13.32'
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("D21").Value = '6.77'
$ws.Range("E21").Value = '  +0.30%  '

$ws.Range("D22").Value = '0.0₃0954'
$ws.Range("E22").Value = '  +0.74%  '

$ws.Range("D23").Value = '281.50'
$ws.Range("E23").Value = '  +0.70%  '

$ws.Range("D24").Value = '72.75'
$ws.Range("E24").Value = '  +1.01%  '

$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").Value = '26.64'
$ws.Range("E26").Value = '  +3.28%  '

$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  -2.10%  '

$ws.Range("E29").Value = '  +1.98%  '

$ws.Range("E30").Value = '  +2.99%  '

$ws.Range("D31").Value = '36.18'
$ws.Range("E31").Value = '  +2.60%  '

$ws.Range("E32").Value = '  +0.58%  '

$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").Value = '19.65'
$ws.Range("E33").Value = '  +0.61%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.46'
$ws.Range("E34").Value = '  +1.89%  '

$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("E36").Value = '  +1.44%  '

$ws.Range("D37").Value = '2.05'
$ws.Range("E37").Value = '  +5.01%  '

$ws.Range("D38").Value = '4.76'
$ws.Range("E38").Value = '  +2.02%  '

$ws.Range("E39").Value = '  +5.76%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.113'
$ws.Range("E40").Value = '  +0.82%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '22.73'
$ws.Range("E41").Value = '  +5.45%  '

$ws.Range("D42").Value = '123.45'
$ws.Range("E42").Value = '  +1.56%  '

$ws.Range("E43").Value = '  +0.47%  '

$ws.Range("E44").Value = '  +3.71%  '

$ws.Range("D45").Value = '3.36'
$ws.Range("E45").Value = '  +5.77%  '

$ws.Range("D46").Value = '2.054.15'
$ws.Range("E46").Value = '  +2.20%  '

$ws.Range("D47").Value = '2.21'
$ws.Range("E47").Value = '  +11.20%  '

$ws.Range("E48").Value = '  +8.93%  '

$ws.Range("E49").Value = '  +0.58%  '

$ws.Range("D50").Value = '5.37'
$ws.Range("E50").Value = '  +3.73%  '

$ws.Range("D51").Value = '81.96'
$ws.Range("E51").Value = '  +2.00%  '

foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
